# Auto-generated edit script for OpenData_Slovakia_Covid_DailyStats
# Commit: "Updated: po 21. 02. 2022"
# - Revises AgTests/AgPosit (columns F/G) figures for several existing rows
# - Appends three new daily rows (716-718) for 2022-02-18 .. 2022-02-20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revisions to existing rows (columns F = AgTests, G = AgPosit) ---
$ws.Cells.Item(617, 6).Value = 39176
$ws.Cells.Item(618, 6).Value = 38098
$ws.Cells.Item(620, 6).Value = 26287
$ws.Cells.Item(621, 6).Value = 56520
$ws.Cells.Item(622, 6).Value = 41699
$ws.Cells.Item(624, 6).Value = 51675
$ws.Cells.Item(625, 6).Value = 44067
$ws.Cells.Item(628, 6).Value = 64877
$ws.Cells.Item(630, 6).Value = 46884
$ws.Cells.Item(631, 6).Value = 42103
$ws.Cells.Item(634, 6).Value = 47155
$ws.Cells.Item(635, 6).Value = 83639
$ws.Cells.Item(638, 6).Value = 37763
$ws.Cells.Item(639, 6).Value = 40739
$ws.Cells.Item(641, 6).Value = 34410
$ws.Cells.Item(642, 6).Value = 67489
$ws.Cells.Item(643, 6).Value = 43478
$ws.Cells.Item(645, 6).Value = 35751
$ws.Cells.Item(646, 6).Value = 36049
$ws.Cells.Item(649, 6).Value = 62676
$ws.Cells.Item(650, 6).Value = 38065
$ws.Cells.Item(652, 6).Value = 35210
$ws.Cells.Item(653, 6).Value = 34159
$ws.Cells.Item(656, 6).Value = 52565
$ws.Cells.Item(657, 6).Value = 34087
$ws.Cells.Item(659, 6).Value = 26391
$ws.Cells.Item(663, 6).Value = 37264
$ws.Cells.Item(666, 6).Value = 23981
$ws.Cells.Item(670, 6).Value = 52650
$ws.Cells.Item(677, 6).Value = 56194
$ws.Cells.Item(680, 6).Value = 28475
$ws.Cells.Item(684, 6).Value = 57257
$ws.Cells.Item(685, 6).Value = 34477
$ws.Cells.Item(686, 6).Value = 34438
$ws.Cells.Item(687, 6).Value = 31482
$ws.Cells.Item(688, 6).Value = 32042
$ws.Cells.Item(691, 6).Value = 62386
$ws.Cells.Item(692, 6).Value = 41604
$ws.Cells.Item(693, 6).Value = 39447
$ws.Cells.Item(694, 6).Value = 37479
$ws.Cells.Item(695, 6).Value = 37111
$ws.Cells.Item(695, 7).Value = 3126
$ws.Cells.Item(696, 6).Value = 17747
$ws.Cells.Item(696, 7).Value = 2207
$ws.Cells.Item(697, 6).Value = 28719
$ws.Cells.Item(697, 7).Value = 3018
$ws.Cells.Item(698, 6).Value = 70032
$ws.Cells.Item(698, 7).Value = 5798
$ws.Cells.Item(699, 6).Value = 43163
$ws.Cells.Item(699, 7).Value = 4278
$ws.Cells.Item(700, 6).Value = 43318
$ws.Cells.Item(700, 7).Value = 4267
$ws.Cells.Item(701, 6).Value = 41432
$ws.Cells.Item(701, 7).Value = 3817
$ws.Cells.Item(702, 6).Value = 36117
$ws.Cells.Item(702, 7).Value = 3892
$ws.Cells.Item(703, 6).Value = 16837
$ws.Cells.Item(703, 7).Value = 2573
$ws.Cells.Item(704, 6).Value = 24687
$ws.Cells.Item(704, 7).Value = 3644
$ws.Cells.Item(705, 6).Value = 55436
$ws.Cells.Item(705, 7).Value = 6237
$ws.Cells.Item(706, 6).Value = 40341
$ws.Cells.Item(706, 7).Value = 4901
$ws.Cells.Item(707, 6).Value = 38303
$ws.Cells.Item(707, 7).Value = 4553
$ws.Cells.Item(708, 6).Value = 35333
$ws.Cells.Item(708, 7).Value = 4118
$ws.Cells.Item(709, 6).Value = 31981
$ws.Cells.Item(709, 7).Value = 3909
$ws.Cells.Item(710, 6).Value = 14493
$ws.Cells.Item(710, 7).Value = 2603
$ws.Cells.Item(711, 6).Value = 22331
$ws.Cells.Item(711, 7).Value = 3782
$ws.Cells.Item(712, 6).Value = 50694
$ws.Cells.Item(712, 7).Value = 6235
$ws.Cells.Item(713, 6).Value = 36520
$ws.Cells.Item(713, 7).Value = 4682
$ws.Cells.Item(714, 6).Value = 30496
$ws.Cells.Item(714, 7).Value = 3627
$ws.Cells.Item(715, 6).Value = 29781
$ws.Cells.Item(715, 7).Value = 3191

# --- New rows appended with latest daily data ---
# Row 716
$ws.Cells.Item(716, 1).Value = 44610
$ws.Cells.Item(716, 2).Value = 1342436
$ws.Cells.Item(716, 3).Value = 32332
$ws.Cells.Item(716, 4).Value = 17933
$ws.Cells.Item(716, 5).Value = 18225
$ws.Cells.Item(716, 6).Value = 23139
$ws.Cells.Item(716, 7).Value = 2855
$ws.Cells.Item(716, 1).NumberFormat = "yyyy-mm-dd"
# Row 717
$ws.Cells.Item(717, 1).Value = 44611
$ws.Cells.Item(717, 2).Value = 1354348
$ws.Cells.Item(717, 3).Value = 22166
$ws.Cells.Item(717, 4).Value = 11912
$ws.Cells.Item(717, 5).Value = 18240
$ws.Cells.Item(717, 6).Value = 9074
$ws.Cells.Item(717, 7).Value = 1584
$ws.Cells.Item(717, 1).NumberFormat = "yyyy-mm-dd"
# Row 718
$ws.Cells.Item(718, 1).Value = 44612
$ws.Cells.Item(718, 2).Value = 1361039
$ws.Cells.Item(718, 3).Value = 12131
$ws.Cells.Item(718, 4).Value = 6691
$ws.Cells.Item(718, 5).Value = 18252
$ws.Cells.Item(718, 6).Value = 8218
$ws.Cells.Item(718, 7).Value = 1473
$ws.Cells.Item(718, 1).NumberFormat = "yyyy-mm-dd"

